# instrument_diagrams.pptx: "Changed pipelined xengine diagram so that N is
# number of ants, and M is accumulation length" + bumped the cached
# datetimeFigureOut field text (28/03/2012 -> 17/04/2012) that lives on the
# slide master and every slide layout.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Cached date field text: 28/03/2012 -> 17/04/2012
#    Present once on the slide master and once on every slide layout, as
#    the literal text of the "Date Placeholder" shape. The placeholder is
#    not always at the same shape index across layouts, so search for the
#    shape whose text matches rather than hard-coding an index.
# ---------------------------------------------------------------------
function Update-DateShapes($container) {
    for ($j = 1; $j -le $container.Shapes.Count; $j++) {
        $shp = $container.Shapes.Item($j)
        if ($shp.HasTextFrame) {
            if ($shp.TextFrame.TextRange.Text -eq "28/03/2012") {
                $shp.TextFrame.TextRange.Text = "17/04/2012"
            }
        }
    }
}

$master = $p.SlideMaster
Update-DateShapes $master
for ($i = 1; $i -le $master.CustomLayouts.Count; $i++) {
    Update-DateShapes $master.CustomLayouts.Item($i)
}

# ---------------------------------------------------------------------
# 2) Pipelined X-Engine diagram (slide 5): swap the roles of N and M, and
#    widen/shift a handful of label text boxes + connectors to match the
#    new (slightly wider) labels.
# ---------------------------------------------------------------------
$s = $p.Slides.Item(5)

# z-N -> z-M (three equivalent textboxes, one per diagram column)
foreach ($idx in 11, 53, 73) {
    $sh = $s.Shapes.Item($idx)
    $sh.TextFrame.TextRange.Characters(2, 2).Text = "-M"
    $sh.Width = 31.068
}

# Elbow connectors that widen to match the wider "z-M" label above them
foreach ($idx in 13, 54, 74) {
    $s.Shapes.Item($idx).Width = 10.10961
}

# Lone "N" -> "M" labels under the delay taps (four equivalent textboxes)
foreach ($idx in 30, 32, 58, 77) {
    $sh = $s.Shapes.Item($idx)
    $sh.TextFrame.TextRange.Text = "M"
    $sh.Width = 21.609
}

# "Cross Tap M/2" -> "Cross Tap N/2"
$crossTap = $s.Shapes.Item(28)
$crossTap.TextFrame.TextRange.Characters(11, 2).Text = "N/"
$crossTap.Left = 512.9824
$crossTap.Width = 85.6039

# Title caption: "  N time sample per antenna, M antennas"
#             -> "  M time sample per antenna, N antennas"
$caption = $s.Shapes.Item(29)
$caption.TextFrame.TextRange.Characters(22, 2).Text = "M "
$caption.TextFrame.TextRange.Characters(49, 1).Text = "N"
$caption.TextFrame.TextRange.Characters(50, 1).Text = " "
